$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Insert-ParagraphXml($xmlFrag) {
    $count = $d.Paragraphs.Count
    $lastPara = $d.Paragraphs.Item($count)
    $endRng = $lastPara.Range.Duplicate
    $endRng.Collapse(0)
    $endRng.InsertParagraphAfter()

    $targetPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $targetRange = $targetPara.Range.Duplicate
    $targetRange.InsertXML($xmlFrag)
}

# Paragraph 1: blank spacer (ListParagraph style, no numbering)
$p1 = '<w:p ' + $wns + '>' +
        '<w:pPr>' +
          '<w:pStyle w:val="ListParagraph"/>' +
          '<w:spacing w:line="276" w:lineRule="auto"/>' +
          '<w:jc w:val="both"/>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>' +
            '<w:color w:val="202122"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
          '</w:rPr>' +
        '</w:pPr>' +
      '</w:p>'
Insert-ParagraphXml($p1)

# Paragraph 2: bold numbered question - "What is Lifting State Up in React?"
$p2 = '<w:p ' + $wns + '>' +
        '<w:pPr>' +
          '<w:pStyle w:val="ListParagraph"/>' +
          '<w:numPr>' +
            '<w:ilvl w:val="0"/>' +
            '<w:numId w:val="1"/>' +
          '</w:numPr>' +
          '<w:spacing w:line="276" w:lineRule="auto"/>' +
          '<w:jc w:val="both"/>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>' +
            '<w:b/>' +
            '<w:bCs/>' +
            '<w:color w:val="202122"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
          '</w:rPr>' +
        '</w:pPr>' +
        '<w:r>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>' +
            '<w:b/>' +
            '<w:bCs/>' +
            '<w:color w:val="202122"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
          '</w:rPr>' +
          '<w:t>What is Lifting State Up in React?</w:t>' +
        '</w:r>' +
      '</w:p>'
Insert-ParagraphXml($p2)

# Paragraph 3: answer text, indented, no list style
$p3 = '<w:p ' + $wns + '>' +
        '<w:pPr>' +
          '<w:spacing w:line="276" w:lineRule="auto"/>' +
          '<w:ind w:left="360"/>' +
          '<w:jc w:val="both"/>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>' +
            '<w:color w:val="202122"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
          '</w:rPr>' +
        '</w:pPr>' +
        '<w:r>' +
          '<w:rPr>' +
            '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/>' +
            '<w:color w:val="202122"/>' +
            '<w:sz w:val="24"/>' +
            '<w:szCs w:val="24"/>' +
            '<w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>' +
          '</w:rPr>' +
          '<w:t>When several components need to share the same changing data then it is recommended to lift the shared state up to their closest common ancestor. That means if two child components share the same data from its parent, then move the state to parent instead of maintaining local state in both of the child components.</w:t>' +
        '</w:r>' +
      '</w:p>'
Insert-ParagraphXml($p3)

Write-Output "Paragraphs now: $($d.Paragraphs.Count)"
